$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right count 5 -> 4, Wrong penalty -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total): Right total 100 -> 80, Wrong penalty -1 -> -2
$ws.Range("B12").Value = 80
$ws.Range("C12").Value = -2

# Update the displayed "score / max" text to reflect corrected totals
$ws.Range("E12").Value = "78 / 112"
